$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the three introductory paragraphs at the top of the
#    document (the "Wearable devices..." / "Through the bracelet..."
#    / "  This device uses..." paragraphs).
# ---------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p3 = $d.Paragraphs.Item(3)
$introRange = $d.Range($p1.Range.Start, $p3.Range.End)
$introRange.Delete()
